$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "328.16"
$cell.Style = $style
$cell = $ws.Range("E2")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.10%"
$cell.Style = $style

$cell = $ws.Range("E3")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.56%"
$cell.Style = $style

$cell = $ws.Range("D4")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.486"
$cell.Style = $style
$cell = $ws.Range("E4")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-1.44%"
$cell.Style = $style

$cell = $ws.Range("D5")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08046"
$cell.Style = $style
$cell = $ws.Range("E5")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.22%"
$cell.Style = $style

$cell = $ws.Range("D6")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.099"
$cell.Style = $style
$cell = $ws.Range("E6")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.69%"
$cell.Style = $style

$cell = $ws.Range("D7")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9524"
$cell.Style = $style
$cell = $ws.Range("E7")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.03%"
$cell.Style = $style

$cell = $ws.Range("D8")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1118"
$cell.Style = $style
$cell = $ws.Range("E8")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-5.46%"
$cell.Style = $style

$cell = $ws.Range("D9")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1878"
$cell.Style = $style
$cell = $ws.Range("E9")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.38%"
$cell.Style = $style

$cell = $ws.Range("D10")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "10.16"
$cell.Style = $style
$cell = $ws.Range("E10")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.83%"
$cell.Style = $style

$cell = $ws.Range("D11")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.09977"
$cell.Style = $style
$cell = $ws.Range("E11")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.21%"
$cell.Style = $style

$cell = $ws.Range("D12")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.04752"
$cell.Style = $style
$cell = $ws.Range("E12")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.94%"
$cell.Style = $style

$cell = $ws.Range("D13")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1059"
$cell.Style = $style
$cell = $ws.Range("E13")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.69%"
$cell.Style = $style

$cell = $ws.Range("D14")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.001272"
$cell.Style = $style
$cell = $ws.Range("E14")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-1.11%"
$cell.Style = $style

$cell = $ws.Range("D15")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.04092"
$cell.Style = $style
$cell = $ws.Range("E15")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-2.60%"
$cell.Style = $style

$cell = $ws.Range("D16")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.006024"
$cell.Style = $style
$cell = $ws.Range("E16")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.96%"
$cell.Style = $style

$cell = $ws.Range("D17")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.373"
$cell.Style = $style
$cell = $ws.Range("E17")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.67%"
$cell.Style = $style

$cell = $ws.Range("D18")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.415"
$cell.Style = $style
$cell = $ws.Range("E18")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.65%"
$cell.Style = $style

$cell = $ws.Range("D19")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.633"
$cell.Style = $style
$cell = $ws.Range("E19")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.28%"
$cell.Style = $style

$cell = $ws.Range("D20")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3284"
$cell.Style = $style
$cell = $ws.Range("E20")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-5.18%"
$cell.Style = $style

$cell = $ws.Range("D21")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1401"
$cell.Style = $style
$cell = $ws.Range("E21")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-1.22%"
$cell.Style = $style

$cell = $ws.Range("D22")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.2494"
$cell.Style = $style
$cell = $ws.Range("E22")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.53%"
$cell.Style = $style

$cell = $ws.Range("D23")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.001311"
$cell.Style = $style
$cell = $ws.Range("E23")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.12%"
$cell.Style = $style

$cell = $ws.Range("D24")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.004345"
$cell.Style = $style
$cell = $ws.Range("E24")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.21%"
$cell.Style = $style

$cell = $ws.Range("D25")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0001253"
$cell.Style = $style
$cell = $ws.Range("E25")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.25%"
$cell.Style = $style

$cell = $ws.Range("D26")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0003745"
$cell.Style = $style
$cell = $ws.Range("E26")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-5.87%"
$cell.Style = $style

$cell = $ws.Range("D38")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.02624"
$cell.Style = $style
$cell = $ws.Range("E38")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-2.16%"
$cell.Style = $style

$cell = $ws.Range("D39")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.05621"
$cell.Style = $style
$cell = $ws.Range("E39")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.47%"
$cell.Style = $style

$cell = $ws.Range("D40")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.007619"
$cell.Style = $style
$cell = $ws.Range("E40")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.70%"
$cell.Style = $style

$cell = $ws.Range("D41")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1401"
$cell.Style = $style
$cell = $ws.Range("E41")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.35%"
$cell.Style = $style

$cell = $ws.Range("D42")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.007379"
$cell.Style = $style
$cell = $ws.Range("E42")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-11.41%"
$cell.Style = $style

$cell = $ws.Range("D43")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.001988"
$cell.Style = $style
$cell = $ws.Range("E43")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-1.42%"
$cell.Style = $style

$cell = $ws.Range("D44")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.008869"
$cell.Style = $style
$cell = $ws.Range("E44")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.47%"
$cell.Style = $style

$cell = $ws.Range("D45")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00007088"
$cell.Style = $style
$cell = $ws.Range("E45")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.28%"
$cell.Style = $style

$cell = $ws.Range("D46")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00000000751"
$cell.Style = $style
$cell = $ws.Range("E46")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.10%"
$cell.Style = $style

$cell = $ws.Range("D47")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0005809"
$cell.Style = $style
$cell = $ws.Range("E47")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "-0.05%"
$cell.Style = $style

$cell = $ws.Range("D48")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.003500"
$cell.Style = $style
$cell = $ws.Range("E48")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "54.06%"
$cell.Style = $style

$cell = $ws.Range("D49")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.003497"
$cell.Style = $style
$cell = $ws.Range("E49")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "16.77%"
$cell.Style = $style

$cell = $ws.Range("D50")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.00002103"
$cell.Style = $style
$cell = $ws.Range("E50")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.10%"
$cell.Style = $style

$cell = $ws.Range("D51")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0002003"
$cell.Style = $style
$cell = $ws.Range("E51")
$style = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.10%"
$cell.Style = $style
